# Insert a new data row at row 31 (shifts existing rows 31-97 down to 32-98)
# and populate it with the new weekly record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(31).Insert()

$ws.Range("A31").Value = 9
$ws.Range("B31").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 44708
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = 100112005
$ws.Range("G31").Value = "Puerro"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 7000
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = 7000
$ws.Range("N31").Value = "$/paquete 20 unidades"
$ws.Range("O31").Value = "Provincia de Chacabuco"
$ws.Range("P31").Value = 350
$ws.Range("Q31").Value = 20
$ws.Range("R31").Value = "Hortaliza"
